$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.757.44'
$ws.Range("E2").Value = '  -3.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.488.66'
$ws.Range("E3").Value = '  -5.88%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.40'
$ws.Range("E5").Value = '  -4.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.61'
$ws.Range("E6").Value = '  -5.02%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.488.31'
$ws.Range("E9").Value = '  -5.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.109'
$ws.Range("E10").Value = '  -7.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("E11").Value = '  -5.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.155'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.359'
$ws.Range("E13").Value = '  -6.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.49'
$ws.Range("E14").Value = '  -7.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.934.32'
$ws.Range("E15").Value = '  -6.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -8.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.589.39'
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.483.69'
$ws.Range("E18").Value = '  -6.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.28'
$ws.Range("E19").Value = '  -7.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.11'
$ws.Range("E20").Value = '  -7.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.23'
$ws.Range("E21").Value = '  -6.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '323.63'
$ws.Range("E22").Value = '  -6.46%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.92'
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.32'
$ws.Range("E25").Value = '  -5.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000102'
$ws.Range("E26").Value = '  -9.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.600.62'
$ws.Range("E27").Value = '  -6.22%  '
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '554.66'
$ws.Range("E28").Value = '  -7.78%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.53'
$ws.Range("E29").Value = '  -5.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.39'
$ws.Range("E32").Value = '  -9.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.153'
$ws.Range("E33").Value = '  -5.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  -7.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  -7.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  -9.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.96'
$ws.Range("E37").Value = '  -9.68%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.386'
$ws.Range("E39").Value = '  -4.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.70'
$ws.Range("E40").Value = '  -5.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '147.19'
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.78'
$ws.Range("E42").Value = '  -6.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.57'
$ws.Range("E44").Value = '  -3.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  -4.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.65'
$ws.Range("E46").Value = '  -9.04%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.25'
$ws.Range("E47").Value = '  -7.84%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.67'
$ws.Range("E48").Value = '  -5.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0546'
$ws.Range("E49").Value = '  -7.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.600'
$ws.Range("E50").Value = '  -5.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0948'
$ws.Range("E51").Value = '  -5.23%  '
